# Owners sheet update: the "system output" API is available now, so we can
# drop the old exporter/exchanger sample rows and record each owner's
# mobile number instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now obsolete sample rows (4444/5555/6666/7777), keeping only
# the header row plus the three remaining owners (1111/2222/3333).
$ws.Range("A5:F8").EntireRow.Delete()

# Add the new "mobileNumber" column header.
$ws.Range("C1").Value = "mobileNumber"

# Row 2 (bankAccountId 1111): still an importer, now with a mobile number.
$ws.Range("B2").Value = "importer"
$ws.Range("C2").Value = "'09059242876"

# Row 3 (bankAccountId 2222): now an exchanger, with a mobile number.
$ws.Range("B3").Value = "exchanger"
$ws.Range("C3").Value = "'09059242876"

# Row 4 (bankAccountId 3333): still an exporter, now with a mobile number.
$ws.Range("B4").Value = "exporter"
$ws.Range("C4").Value = "'09059242876"

# Match the author's last selected cell.
$ws.Range("D9").Select() | Out-Null
